$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.237.03'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.427.34'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.50'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.49'
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.428.32'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.559'
$ws.Range('E9').Value = '  -7.95%  '
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.427'
$ws.Range('E12').Value = '  -3.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.024.20'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.26'
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  -6.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.296.62'
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.456.08'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('E19').Value = '  -3.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.66'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.83'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.89'
$ws.Range('E22').Value = '  -1.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.60'
$ws.Range('E24').Value = '  -0.98%  '
$ws.Range('E25').Value = '  -5.14%  '
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.63'
$ws.Range('E27').Value = '  -4.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.177'
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('E31').Value = '  -4.05%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.04'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.11'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('E35').Value = '  -4.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.04'
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.858'
$ws.Range('E37').Value = '  +11.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.83'
$ws.Range('E38').Value = '  -4.38%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0734'
$ws.Range('E39').Value = '  -2.99%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.805.79'
$ws.Range('E40').Value = '  -3.94%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '25.94'
$ws.Range('E41').Value = '  -2.17%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.69'
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.11'
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.50'
$ws.Range('E44').Value = '  -3.62%  '
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('E46').Value = '  -3.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '344.42'
$ws.Range('E47').Value = '  +7.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  +8.44%  '
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.35'
$ws.Range('E50').Value = '  -3.01%  '
$ws.Range('E51').Value = '  -4.14%  '
